# Apply cryptos list update (price + volume refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.465.19"
$ws.Range("E2").Value = "  +3.52%  "

$ws.Range("D3").Value = "3.067.05"
$ws.Range("E3").Value = "  +2.83%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.96"
$ws.Range("E5").Value = "  +3.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.21"
$ws.Range("E6").Value = "  +7.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.063.75"
$ws.Range("E8").Value = "  +2.87%  "

$ws.Range("E9").Value = "  +1.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.52"
$ws.Range("E10").Value = "  +6.68%  "

$ws.Range("E11").Value = "  +2.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  +3.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.83"
$ws.Range("E14").Value = "  +3.88%  "

$ws.Range("D15").Value = "3.567.36"
$ws.Range("E15").Value = "  +2.91%  "

$ws.Range("D16").Value = "63.446.88"
$ws.Range("E16").Value = "  +3.41%  "

$ws.Range("D17").Value = "3.071.88"
$ws.Range("E17").Value = "  +2.95%  "

$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("E19").Value = "  +3.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.40"
$ws.Range("E20").Value = "  +4.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.89"
$ws.Range("E21").Value = "  +5.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.677"
$ws.Range("E22").Value = "  +1.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.28"
$ws.Range("E23").Value = "  +5.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.90"
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.84"
$ws.Range("E25").Value = "  +8.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  +4.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("E28").Value = "  +3.19%  "

$ws.Range("E29").Value = "  +8.13%  "

$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.15"
$ws.Range("E31").Value = "  +2.85%  "

$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("E33").Value = "  +8.23%  "

$ws.Range("E34").Value = "  +5.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.36"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.01"
$ws.Range("E36").Value = "  +2.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "467.03"
$ws.Range("E37").Value = "  +2.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0824"
$ws.Range("E38").Value = "  +5.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0396"
$ws.Range("E39").Value = "  +3.96%  "

$ws.Range("D40").Value = "3.016.80"
$ws.Range("E40").Value = "  -4.11%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.20"
$ws.Range("E42").Value = "  +1.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  +6.10%  "

$ws.Range("E44").Value = "  +5.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.255"
$ws.Range("E45").Value = "  +6.13%  "

$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.04"
$ws.Range("E47").Value = "  +2.87%  "

$ws.Range("E48").Value = "  +2.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.60"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("E50").Value = "  +3.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("E51").Value = "  +4.68%  "
